$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap firstName (column B) and lastName (column C) for rows 12 through 51
for ($r = 12; $r -le 51; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value()
    $cVal = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($r, 2).Value = $cVal
    $ws.Cells.Item($r, 3).Value = $bVal
}

# Update the view: scroll so row 25 is at the top-left, and select L38
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("L38").Select()
